$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 7782
$ws.Range("F3").Value = 7782
$ws.Range("F4").Value = 15
$ws.Range("F5").Value = 7937
$ws.Range("F8").Value = 39
$ws.Range("F9").Value = 6804
$ws.Range("F10").Value = 3416
$ws.Range("F12").Value = 3751
$ws.Range("F13").Value = 52
$ws.Range("F14").Value = 58
$ws.Range("F15").Value = 52
$ws.Range("F16").Value = 78
$ws.Range("F17").Value = 87
$ws.Range("F18").Value = 482
$ws.Range("F20").Value = 64
$ws.Range("F24").Value = 340
$ws.Range("F25").Value = 3912
$ws.Range("F28").Value = 1014
$ws.Range("F29").Value = 500
$ws.Range("F30").Value = 1531
$ws.Range("F32").Value = 70
$ws.Range("F33").Value = 2813
$ws.Range("F34").Value = 1976
$ws.Range("F35").Value = 42
$ws.Range("F36").Value = 59
$ws.Range("F39").Value = 3826
$ws.Range("F41").Value = 289
$ws.Range("F45").Value = 18
$ws.Range("F46").Value = 1473
$ws.Range("F48").Value = 571
$ws.Range("F49").Value = 663

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 28
$ws.Range("F10").Value = 24
$ws.Range("F17").Value = 243
$ws.Range("F18").Value = 2

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 7782
$ws.Range("F6").Value = 7782
$ws.Range("F7").Value = 7937
$ws.Range("F8").Value = 39
$ws.Range("F9").Value = 6804
$ws.Range("F10").Value = 3416
$ws.Range("F11").Value = 3751
$ws.Range("F12").Value = 58
$ws.Range("F13").Value = 52
$ws.Range("F14").Value = 78
$ws.Range("F15").Value = 87
$ws.Range("F16").Value = 482
$ws.Range("F18").Value = 64
$ws.Range("F21").Value = 340
$ws.Range("F22").Value = 3912
$ws.Range("F23").Value = 24
$ws.Range("F27").Value = 1014
$ws.Range("F28").Value = 500
$ws.Range("F29").Value = 1531
$ws.Range("F31").Value = 70
$ws.Range("F32").Value = 2813
$ws.Range("F33").Value = 1976
$ws.Range("F34").Value = 42
$ws.Range("F35").Value = 59
$ws.Range("F39").Value = 3826
$ws.Range("F41").Value = 289
$ws.Range("F45").Value = 243
$ws.Range("F46").Value = 1473
$ws.Range("F48").Value = 2
$ws.Range("F49").Value = 571
$ws.Range("F50").Value = 663
